$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 76925256
$ws.Range("I18").Value = 90911400
$ws.Range("J18").Value = 1446
$ws.Range("K18").Value = 90911400
$ws.Range("L18").Value = 1446
$ws.Range("M18").Value = -90911116
$ws.Range("N18").Value = -2014
$ws.Range("H40").Value = 3953.8462
$ws.Range("J40").Value = 5362.5
$ws.Range("L40").Value = 5362.5
$ws.Range("N40").Value = -5712.5
$ws.Range("H70").Value = 2611.2144
$ws.Range("I70").Value = 844
$ws.Range("J70").Value = 3318.1
$ws.Range("K70").Value = 2532
$ws.Range("L70").Value = 9954.299999999999
$ws.Range("M70").Value = -2262
$ws.Range("N70").Value = -10494.3
$ws.Range("H73").Value = 2611.2144
$ws.Range("I73").Value = 844
$ws.Range("J73").Value = 3318.1
$ws.Range("K73").Value = 2532
$ws.Range("L73").Value = 9954.299999999999
$ws.Range("M73").Value = -1596
$ws.Range("N73").Value = -11826.3
$ws.Range("H86").Value = 2044.375
$ws.Range("I86").Value = 2044.375
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2044.375
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -921.375
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2044.375
$ws.Range("I89").Value = 2044.375
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10221.875
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4605.875
$ws.Range("N89").ClearContents()
$ws.Range("H140").Value = 69215.60000000001
$ws.Range("J140").Value = 69215.60000000001
$ws.Range("L140").Value = 69215.60000000001
$ws.Range("N140").Value = -79575.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 13367.462
$ws.Range("J31").Value = 18497.5
$ws.Range("L31").Value = 18497.5
$ws.Range("N31").Value = -19085.5
$ws.Range("H32").Value = 21398.82
$ws.Range("J32").Value = 36720.586
$ws.Range("L32").Value = 36720.586
$ws.Range("N32").Value = -37294.586
$ws.Range("H35").Value = 3999
$ws.Range("I35").Value = 3999
$ws.Range("K35").Value = 3999
$ws.Range("M35").Value = -3593
$ws.Range("H61").Value = 5551.607
$ws.Range("I61").Value = 4909.8945
$ws.Range("K61").Value = 4909.8945
$ws.Range("M61").Value = -4697.8945
$ws.Range("H88").Value = 2143.125
$ws.Range("J88").Value = 2469.6
$ws.Range("L88").Value = 2469.6
$ws.Range("N88").Value = -3281.6
$ws.Range("H91").Value = 2143.125
$ws.Range("J91").Value = 2469.6
$ws.Range("L91").Value = 2469.6
$ws.Range("N91").Value = -5277.6
$ws.Range("H136").Value = 5551.607
$ws.Range("I136").Value = 4909.8945
$ws.Range("K136").Value = 14729.6835
$ws.Range("M136").Value = -12179.6835

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 11731
$ws.Range("I75").Value = 6219.3335
$ws.Range("J75").Value = 19998.5
$ws.Range("K75").Value = 6219.3335
$ws.Range("L75").Value = 19998.5
$ws.Range("M75").Value = -5283.3335
$ws.Range("N75").Value = -21870.5
$ws.Range("H78").Value = 11731
$ws.Range("I78").Value = 6219.3335
$ws.Range("J78").Value = 19998.5
$ws.Range("K78").Value = 18658.0005
$ws.Range("L78").Value = 59995.5
$ws.Range("M78").Value = -13978.0005
$ws.Range("N78").Value = -69355.5
$ws.Range("H134").Value = 2900.6191
$ws.Range("I134").Value = 1675.4193
$ws.Range("K134").Value = 5026.257900000001
$ws.Range("M134").Value = -2491.257900000001
$ws.Range("H138").Value = 94606.664
$ws.Range("J138").Value = 94606.664
$ws.Range("L138").Value = 94606.664
$ws.Range("N138").Value = -104886.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4123.304
$ws.Range("I31").Value = 3102.4
$ws.Range("J31").Value = 4540
$ws.Range("K31").Value = 3102.4
$ws.Range("L31").Value = 4540
$ws.Range("M31").Value = -2807.4
$ws.Range("N31").Value = -5130
$ws.Range("H34").Value = 4123.304
$ws.Range("I34").Value = 3102.4
$ws.Range("J34").Value = 4540
$ws.Range("K34").Value = 3102.4
$ws.Range("L34").Value = 4540
$ws.Range("M34").Value = -2900.4
$ws.Range("N34").Value = -4944
$ws.Range("H62").Value = 81757.78999999999
$ws.Range("I62").Value = 176817.33
$ws.Range("K62").Value = 176817.33
$ws.Range("M62").Value = -176193.33
$ws.Range("H65").Value = 81757.78999999999
$ws.Range("I65").Value = 176817.33
$ws.Range("K65").Value = 884086.6499999999
$ws.Range("M65").Value = -880966.6499999999
$ws.Range("H138").Value = 35903.617
$ws.Range("J138").Value = 35903.617
$ws.Range("L138").Value = 35903.617
$ws.Range("N138").Value = -46183.617

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 505179.2
$ws.Range("I4").Value = 611588
$ws.Range("K4").Value = 1834764
$ws.Range("M4").Value = -1834652
$ws.Range("H137").Value = 670591.4
$ws.Range("I137").Value = 837215.3
$ws.Range("J137").Value = 4095.6667
$ws.Range("K137").Value = 2511645.9
$ws.Range("L137").Value = 12287.0001
$ws.Range("M137").Value = -2506545.9
$ws.Range("N137").Value = -22487.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 19882
$ws.Range("I80").Value = 1800
$ws.Range("J80").Value = 23498.4
$ws.Range("K80").Value = 1800
$ws.Range("L80").Value = 23498.4
$ws.Range("M80").Value = -802
$ws.Range("N80").Value = -25494.4
$ws.Range("H83").Value = 19882
$ws.Range("I83").Value = 1800
$ws.Range("J83").Value = 23498.4
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 117492
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -127476
$ws.Range("H135").Value = 99835.2
$ws.Range("J135").Value = 99835.2
$ws.Range("L135").Value = 99835.2
$ws.Range("N135").Value = -109975.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 9311.727999999999
$ws.Range("I53").Value = 13983
$ws.Range("J53").Value = 8273.666999999999
$ws.Range("K53").Value = 13983
$ws.Range("L53").Value = 8273.666999999999
$ws.Range("M53").Value = -13465
$ws.Range("N53").Value = -9309.666999999999
$ws.Range("H82").Value = 166669500
$ws.Range("I82").Value = 3400
$ws.Range("K82").Value = 3400
$ws.Range("M82").Value = -3039
$ws.Range("H85").Value = 166669500
$ws.Range("I85").Value = 3400
$ws.Range("K85").Value = 3400
$ws.Range("M85").Value = -2152

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 40716.168
$ws.Range("J103").Value = 40716.168
$ws.Range("L103").Value = 40716.168
$ws.Range("N103").Value = -43060.168
$ws.Range("H139").Value = 93666.336
$ws.Range("I139").Value = 50000
$ws.Range("J139").Value = 102399.6
$ws.Range("K139").Value = 50000
$ws.Range("L139").Value = 102399.6
$ws.Range("M139").Value = -44860
$ws.Range("N139").Value = -112679.6
$ws.Range("H141").Value = 109950.4
$ws.Range("J141").Value = 109950.4
$ws.Range("L141").Value = 109950.4
$ws.Range("N141").Value = -120310.4
